# Generate Report for Handoff
#
# A newer handoff run (file "afe5f849-...") supersedes the previous one
# (file "3a009bb0-..."): the tracked source-file name/hash is updated
# everywhere it is referenced, the handoff timestamps are refreshed, the
# now-stale "Latest Target File" info on the two locale sheets is cleared
# (along with its hyperlink), and those two columns are narrowed to their
# new, tighter auto-fit widths.

$wb = $excel.ActiveWorkbook

$oldGuid = "3a009bb0-b78e-4caf-9503-e1452583d4ef"
$newGuid = "afe5f849-f94e-4a21-bc87-42d5d896504b"
$newHash = "9d69580233f4b663cdac7cdcfe0e6870c99eb70f"

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("G2").Value = "2016-08-29 19:05:39"

$ws1.Range("B2").Value = "e2e\$newGuid.md"
foreach ($h in $ws1.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') {
        $h.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "$newGuid.md"
foreach ($h in $ws2.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    }
}

$ws2.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-29 19:05:34"
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

# "Latest Target File" (I2) is no longer tracked: its hyperlink goes away
# and both it and "Latest Handback File" (J2) are blanked out. The engine
# only supports removing hyperlinks sheet-wide, so capture the A2 link's
# target first, drop every link on the sheet, then restore just the A2 one.
$a2Addr2 = ""
foreach ($h in $ws2.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $a2Addr2 = $h.Address
    }
}
$ws2.Range("A1:P10").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $a2Addr2, "", "", "$newGuid.md")

$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = ""
$ws2.Range("J2").Style = "Normal"

$ws2.Columns.Item(9).ColumnWidth = 17.86
$ws2.Columns.Item(10).ColumnWidth = 20.86

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3.Range("A2").Value = "$newGuid.md"
foreach ($h in $ws3.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    }
}

$ws3.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-29 19:05:39"
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$a2Addr3 = ""
foreach ($h in $ws3.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $a2Addr3 = $h.Address
    }
}
$ws3.Range("A1:P10").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $a2Addr3, "", "", "$newGuid.md")

$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = ""
$ws3.Range("J2").Style = "Normal"

$ws3.Columns.Item(9).ColumnWidth = 17.86
$ws3.Columns.Item(10).ColumnWidth = 20.86
